# Update Main Asset Number (column D) and Asset Subnumber (column E)
# values on the "Data" sheet for the MPA test automation upload file.
#
# Rows 6,7,8,9,10,16,17    -> D (Main Asset Number) increases by 18
# Rows 11,12,13,14,15,18,19 -> E (Asset Subnumber) increases by 18
# Rows 20,22,24,26          -> D increases by 18
# Rows 21,23,25,27          -> E increases by 18

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Column D (Main Asset Number) updates: 60000408 -> 60000426, etc.
$dRows = @(6, 7, 8, 9, 10, 16, 17)
foreach ($r in $dRows) {
    $ws.Range("D$r").Value = 60000426
}

$ws.Range("D20").Value = 60000427
$ws.Range("D22").Value = 60000428
$ws.Range("D24").Value = 60000429
$ws.Range("D26").Value = 60000430

# Column E (Asset Subnumber) updates: 284 -> 291, etc.
$eRows = @(11, 12, 13, 14, 15, 18, 19)
foreach ($r in $eRows) {
    $ws.Range("E$r").Value = 291
}

$ws.Range("E21").Value = 292
$ws.Range("E23").Value = 293
$ws.Range("E25").Value = 294
$ws.Range("E27").Value = 295
